$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 42
$ws.Range("H42").Value = 184.25
$ws.Range("I42").Value = 64
$ws.Range("J42").Value = 304.5
$ws.Range("K42").Value = 192
$ws.Range("L42").Value = 913.5
$ws.Range("M42").Value = 38
$ws.Range("N42").Value = -1373.5

# Row 75
$ws.Range("H75").Value = 109833.336
$ws.Range("I75").Value = 30000
$ws.Range("J75").Value = 149750
$ws.Range("K75").Value = 30000
$ws.Range("L75").Value = 149750
$ws.Range("M75").Value = -29064
$ws.Range("N75").Value = -151622

# Row 78
$ws.Range("H78").Value = 109833.336
$ws.Range("I78").Value = 30000
$ws.Range("J78").Value = 149750
$ws.Range("K78").Value = 90000
$ws.Range("L78").Value = 449250
$ws.Range("M78").Value = -85320
$ws.Range("N78").Value = -458610

# Row 93
$ws.Range("H93").Value = 29020
$ws.Range("J93").Value = 29020
$ws.Range("L93").Value = 29020
$ws.Range("N93").Value = -34012


$ws = $wb.Worksheets.Item("ARM")
# Row 94
$ws.Range("H94").Value = 29500
$ws.Range("J94").Value = 29500
$ws.Range("L94").Value = 29500
$ws.Range("N94").Value = -31302

# Row 124
$ws.Range("H124").Value = 54939.668
$ws.Range("J124").Value = 54939.668
$ws.Range("L124").Value = 54939.668
$ws.Range("N124").Value = -64759.668

# Row 125
$ws.Range("H125").Value = 32992.715
$ws.Range("J125").Value = 32992.715
$ws.Range("L125").Value = 32992.715
$ws.Range("N125").Value = -42832.715


$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 50000
$ws.Range("I2").Value = 50000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 50000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -49887
$ws.Range("N2").ClearContents()

# Row 28
$ws.Range("H28").Value = 29771.5
$ws.Range("J28").Value = 29771.5
$ws.Range("L28").Value = 29771.5
$ws.Range("N28").Value = -30261.5

# Row 93
$ws.Range("H93").Value = 17985.143
$ws.Range("I93").Value = 15379.2
$ws.Range("J93").Value = 24500
$ws.Range("K93").Value = 15379.2
$ws.Range("L93").Value = 24500
$ws.Range("M93").Value = -13507.2
$ws.Range("N93").Value = -28244


$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 12144.76
$ws.Range("I4").Value = 97.61539
$ws.Range("J4").Value = 25195.834
$ws.Range("K4").Value = 292.84617
$ws.Range("L4").Value = 75587.50199999999
$ws.Range("M4").Value = -180.84617
$ws.Range("N4").Value = -75811.50199999999

# Row 7
$ws.Range("H7").Value = 86.5
$ws.Range("I7").Value = 86.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 259.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -147.5
$ws.Range("N7").ClearContents()

# Row 12
$ws.Range("H12").Value = 104.77778
$ws.Range("J12").Value = 90.40000000000001
$ws.Range("L12").Value = 271.2
$ws.Range("N12").Value = -617.2

# Row 13
$ws.Range("H13").Value = 323.9
$ws.Range("I13").Value = 165.57143
$ws.Range("J13").Value = 693.3333
$ws.Range("K13").Value = 496.71429
$ws.Range("L13").Value = 2079.9999
$ws.Range("M13").Value = -328.71429
$ws.Range("N13").Value = -2415.9999

# Row 17
$ws.Range("H17").Value = 1236
$ws.Range("I17").Value = 1793.3334
$ws.Range("J17").Value = 400
$ws.Range("K17").Value = 5380.0002
$ws.Range("L17").Value = 1200
$ws.Range("M17").Value = -5211.0002
$ws.Range("N17").Value = -1538

# Row 23
$ws.Range("H23").Value = 716.7
$ws.Range("I23").Value = 1219.1111
$ws.Range("J23").Value = 305.63635
$ws.Range("K23").Value = 3657.3333
$ws.Range("L23").Value = 916.90905
$ws.Range("M23").Value = -3422.3333
$ws.Range("N23").Value = -1386.90905

# Row 24
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

# Row 25
$ws.Range("H25").Value = 500
$ws.Range("I25").Value = 500
$ws.Range("K25").Value = 1500
$ws.Range("M25").Value = -1331

# Row 30
$ws.Range("H30").Value = 500
$ws.Range("I30").Value = 500
$ws.Range("K30").Value = 1500
$ws.Range("M30").Value = -1398

# Row 41
$ws.Range("H41").Value = 3000
$ws.Range("J41").Value = 3000
$ws.Range("L41").Value = 9000
$ws.Range("N41").Value = -9676

# Row 42
$ws.Range("H42").Value = 1522.2222
$ws.Range("I42").Value = 1500
$ws.Range("J42").Value = 1550
$ws.Range("K42").Value = 4500
$ws.Range("L42").Value = 4650
$ws.Range("M42").Value = -3966
$ws.Range("N42").Value = -5718


$ws = $wb.Worksheets.Item("LTW")
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# Row 76
$ws.Range("H76").Value = 35000
$ws.Range("J76").Value = 35000
$ws.Range("L76").Value = 35000
$ws.Range("N76").Value = -35676

# Row 79
$ws.Range("H79").Value = 35000
$ws.Range("J79").Value = 35000
$ws.Range("L79").Value = 35000
$ws.Range("N79").Value = -37340

# Row 94
$ws.Range("H94").Value = 48995
$ws.Range("J94").Value = 48995
$ws.Range("L94").Value = 48995
$ws.Range("N94").Value = -50347


$ws = $wb.Worksheets.Item("WVR")
# Row 63
$ws.Range("H63").Value = 32412.857
$ws.Range("J63").Value = 32412.857
$ws.Range("L63").Value = 32412.857
$ws.Range("N63").Value = -33660.857

# Row 66
$ws.Range("H66").Value = 32412.857
$ws.Range("J66").Value = 32412.857
$ws.Range("L66").Value = 97238.571
$ws.Range("N66").Value = -103478.571

# Row 92
$ws.Range("H92").Value = 29500
$ws.Range("J92").Value = 29500
$ws.Range("L92").Value = 29500
$ws.Range("N92").Value = -34492

# Row 104
$ws.Range("H104").Value = 30000
$ws.Range("J104").Value = 30000
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -36988

# Row 136
$ws.Range("H136").Value = 2341.476
$ws.Range("I136").Value = 892.4400000000001
$ws.Range("J136").Value = 7914.6924
$ws.Range("K136").Value = 2677.32
$ws.Range("L136").Value = 23744.0772
$ws.Range("M136").Value = -127.3200000000002
$ws.Range("N136").Value = -28844.0772

